$d = $word.ActiveDocument

# Move to the very end of the document body and append:
#   1) a blank paragraph
#   2) a new paragraph containing the INSERT INTO statement, as two runs
#      (the statement text, then a separate run holding just the
#      trailing semicolon)
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.InsertParagraphAfter()

$last = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$last.InsertAfter('INSERT INTO `user` VALUES (1,''admin'',NULL,NULL,''admin@gmail.com'',NULL,''admin'',''admin'',NULL,NULL,NULL,NULL,0);')

# Split the trailing ";" into its own run (matches the source document,
# which has the statement text and the semicolon in separate runs).
$full = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$semi = $d.Range($full.End - 2, $full.End - 1)
$semi.Bold = 1
$semi.Bold = 0
